# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H-N across several rows
# in sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 316.45456
$ws.Range("I5").Value = 341.375
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 341.375
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -226.375
$ws.Range("N5").Value = -480

$ws.Range("H132").Value = 1885720
$ws.Range("I132").Value = 1093.0435
$ws.Range("J132").Value = 16334527
$ws.Range("K132").Value = 3279.1305
$ws.Range("L132").Value = 49003581
$ws.Range("M132").Value = -749.1305000000002
$ws.Range("N132").Value = -49008641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L6").ClearContents()
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -327
$ws.Range("N6").Value = 0

$ws.Range("H61").Value = 22267812
$ws.Range("I61").Value = 27806552
$ws.Range("J61").Value = 112853.336
$ws.Range("K61").Value = 27806552
$ws.Range("L61").Value = 112853.336
$ws.Range("M61").Value = -27806340
$ws.Range("N61").Value = -113277.336

$ws.Range("H74").Value = 4583590
$ws.Range("I74").Value = 6436984
$ws.Range("J74").Value = 65940.94
$ws.Range("K74").Value = 6436984
$ws.Range("L74").Value = 65940.94
$ws.Range("M74").Value = -6436110
$ws.Range("N74").Value = -67688.94

$ws.Range("H77").Value = 4583590
$ws.Range("I77").Value = 6436984
$ws.Range("J77").Value = 65940.94
$ws.Range("K77").Value = 32184920
$ws.Range("L77").Value = 329704.7
$ws.Range("M77").Value = -32180552
$ws.Range("N77").Value = -338440.7

$ws.Range("H122").Value = 2925581
$ws.Range("I122").Value = 1578.3667
$ws.Range("J122").Value = 13890591
$ws.Range("K122").Value = 4735.1001
$ws.Range("L122").Value = 41671773
$ws.Range("M122").Value = -2285.1001
$ws.Range("N122").Value = -41676673

$ws.Range("H132").Value = 39545.207
$ws.Range("I132").Value = 29126.111
$ws.Range("J132").Value = 61609.176
$ws.Range("K132").Value = 87378.333
$ws.Range("L132").Value = 184827.528
$ws.Range("M132").Value = -84848.333
$ws.Range("N132").Value = -189887.528

$ws.Range("H136").Value = 22267812
$ws.Range("I136").Value = 27806552
$ws.Range("J136").Value = 112853.336
$ws.Range("K136").Value = 83419656
$ws.Range("L136").Value = 338560.008
$ws.Range("M136").Value = -83417106
$ws.Range("N136").Value = -343660.008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 25000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 25000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774

$ws.Range("H134").Value = 2152
$ws.Range("I134").Value = 1431.4814
$ws.Range("J134").Value = 6042.8
$ws.Range("K134").Value = 4294.4442
$ws.Range("L134").Value = 18128.4
$ws.Range("M134").Value = -1759.4442
$ws.Range("N134").Value = -23198.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1767
$ws.Range("I31").Value = 1279.2693
$ws.Range("J31").Value = 3176
$ws.Range("K31").Value = 1279.2693
$ws.Range("L31").Value = 3176
$ws.Range("M31").Value = -984.2692999999999
$ws.Range("N31").Value = -3766

$ws.Range("H34").Value = 1767
$ws.Range("I34").Value = 1279.2693
$ws.Range("J34").Value = 3176
$ws.Range("K34").Value = 1279.2693
$ws.Range("L34").Value = 3176
$ws.Range("M34").Value = -1077.2693
$ws.Range("N34").Value = -3580

$ws.Range("H58").Value = 17242830
$ws.Range("I58").Value = 25642054
$ws.Range("J58").Value = 2316.6843
$ws.Range("K58").Value = 25642054
$ws.Range("L58").Value = 2316.6843
$ws.Range("M58").Value = -25641851
$ws.Range("N58").Value = -2722.6843

$ws.Range("H132").Value = 27209.725
$ws.Range("I132").Value = 1955.3846
$ws.Range("J132").Value = 74110.64
$ws.Range("K132").Value = 5866.1538
$ws.Range("L132").Value = 222331.92
$ws.Range("M132").Value = -3336.1538
$ws.Range("N132").Value = -227391.92

$ws.Range("H134").Value = 22422.389
$ws.Range("I134").Value = 1552.8049
$ws.Range("J134").Value = 88241.84
$ws.Range("K134").Value = 4658.4147
$ws.Range("L134").Value = 264725.52
$ws.Range("M134").Value = -2123.4147
$ws.Range("N134").Value = -269795.52

$ws.Range("H136").Value = 17242830
$ws.Range("I136").Value = 25642054
$ws.Range("J136").Value = 2316.6843
$ws.Range("K136").Value = 76926162
$ws.Range("L136").Value = 6950.0529
$ws.Range("M136").Value = -76923612
$ws.Range("N136").Value = -12050.0529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 266.77777
$ws.Range("I33").Value = 133.66667
$ws.Range("J33").Value = 333.33334
$ws.Range("K33").Value = 802.0000200000001
$ws.Range("L33").Value = 2000.00004
$ws.Range("M33").Value = -519.0000200000001
$ws.Range("N33").Value = -2566.00004

$ws.Range("H44").Value = 433.33334
$ws.Range("I44").Value = 400
$ws.Range("J44").Value = 435
$ws.Range("K44").Value = 1200
$ws.Range("L44").Value = 1305
$ws.Range("M44").Value = -802
$ws.Range("N44").Value = -2101

$ws.Range("H58").Value = 668.3333
$ws.Range("I58").Value = 552.5
$ws.Range("J58").Value = 900
$ws.Range("K58").Value = 1657.5
$ws.Range("L58").Value = 2700
$ws.Range("M58").Value = -1529.5
$ws.Range("N58").Value = -2956

$ws.Range("H64").Value = 3564.8147
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 3625
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 10875
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -11415

$ws.Range("H67").Value = 3564.8147
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 3625
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 10875
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -12747

$ws.Range("H68").Value = 435828.3
$ws.Range("I68").Value = 955.1724
$ws.Range("J68").Value = 1177670.8
$ws.Range("K68").Value = 2865.5172
$ws.Range("L68").Value = 3533012.4
$ws.Range("M68").Value = -2054.5172
$ws.Range("N68").Value = -3534634.4

$ws.Range("H70").Value = 3275
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 3600
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 10800
$ws.Range("M70").Value = -2685
$ws.Range("N70").Value = -11430

$ws.Range("H71").Value = 435828.3
$ws.Range("I71").Value = 955.1724
$ws.Range("J71").Value = 1177670.8
$ws.Range("K71").Value = 8596.551600000001
$ws.Range("L71").Value = 10599037.2
$ws.Range("M71").Value = -4540.551600000001
$ws.Range("N71").Value = -10607149.2

$ws.Range("H73").Value = 3275
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 3600
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 10800
$ws.Range("M73").Value = -1908
$ws.Range("N73").Value = -12984

$ws.Range("H76").Value = 3700
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3700
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 11100
$ws.Range("N76").Value = -11866

$ws.Range("H79").Value = 3700
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3700
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 11100
$ws.Range("N79").Value = -13752

$ws.Range("H92").Value = 975.7143
$ws.Range("I92").Value = 1250
$ws.Range("J92").Value = 930
$ws.Range("K92").Value = 3750
$ws.Range("L92").Value = 2790
$ws.Range("M92").Value = -2502
$ws.Range("N92").Value = -5286

$ws.Range("L94").ClearContents()
$ws.Range("H94").Value = 3502.2727
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3502.2727
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = 10506.8181
$ws.Range("N94").Value = -11858.8181

$ws.Range("H100").Value = 2438.5
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 2514.2104
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 7542.6312
$ws.Range("M100").Value = -2189
$ws.Range("N100").Value = -9164.6312

$ws.Range("H103").Value = 2396.1924
$ws.Range("I103").Value = 447.3
$ws.Range("J103").Value = 3614.25
$ws.Range("K103").Value = 1341.9
$ws.Range("L103").Value = 10842.75
$ws.Range("M103").Value = -462.9000000000001
$ws.Range("N103").Value = -12600.75

$ws.Range("H106").Value = 3095.625
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 3095.625
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 9286.875
$ws.Range("N106").Value = -11178.875

$ws.Range("H109").Value = 2747.0476
$ws.Range("I109").Value = 2012.5714
$ws.Range("J109").Value = 3114.2856
$ws.Range("K109").Value = 6037.7142
$ws.Range("L109").Value = 9342.856800000001
$ws.Range("M109").Value = -4997.7142
$ws.Range("N109").Value = -11422.8568

$ws.Range("H112").Value = 13336371
$ws.Range("I112").Value = 2047.4445
$ws.Range("J112").Value = 20836928
$ws.Range("K112").Value = 6142.333500000001
$ws.Range("L112").Value = 62510784
$ws.Range("M112").Value = -5034.333500000001
$ws.Range("N112").Value = -62513000

$ws.Range("H131").Value = 15712.442
$ws.Range("I131").Value = 865
$ws.Range("J131").Value = 17332.164
$ws.Range("K131").Value = 2595
$ws.Range("L131").Value = 51996.492
$ws.Range("M131").Value = 2445
$ws.Range("N131").Value = -62076.492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2577.1667
$ws.Range("I126").Value = 2268.6667
$ws.Range("J126").Value = 2885.6667
$ws.Range("K126").Value = 6806.000100000001
$ws.Range("L126").Value = 8657.000100000001
$ws.Range("M126").Value = -4336.000100000001
$ws.Range("N126").Value = -13597.0001

$ws.Range("H132").Value = 88948.83
$ws.Range("I132").Value = 126141
$ws.Range("J132").Value = 69113
$ws.Range("K132").Value = 378423
$ws.Range("L132").Value = 207339
$ws.Range("M132").Value = -375893
$ws.Range("N132").Value = -212399

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2144.818
$ws.Range("I7").Value = 2027.1428
$ws.Range("J7").Value = 2350.75
$ws.Range("K7").Value = 2027.1428
$ws.Range("L7").Value = 2350.75
$ws.Range("M7").Value = -1915.1428
$ws.Range("N7").Value = -2574.75

$ws.Range("H126").Value = 2144.818
$ws.Range("I126").Value = 2027.1428
$ws.Range("J126").Value = 2350.75
$ws.Range("K126").Value = 6081.428400000001
$ws.Range("L126").Value = 7052.25
$ws.Range("M126").Value = -3611.428400000001
$ws.Range("N126").Value = -11992.25

$ws.Range("H132").Value = 16806.322
$ws.Range("I132").Value = 1127.804
$ws.Range("J132").Value = 73920.92999999999
$ws.Range("K132").Value = 3383.412
$ws.Range("L132").Value = 221762.79
$ws.Range("M132").Value = -853.4120000000003
$ws.Range("N132").Value = -226822.79

$ws.Range("H136").Value = 101367.1
$ws.Range("I136").Value = 67815.87
$ws.Range("J136").Value = 202020.8
$ws.Range("K136").Value = 203447.61
$ws.Range("L136").Value = 606062.3999999999
$ws.Range("M136").Value = -200897.61
$ws.Range("N136").Value = -611162.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L53").ClearContents()
$ws.Range("H53").Value = 35000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 35000
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = 35000
$ws.Range("N53").Value = -36214

$ws.Range("H126").Value = 1338.7333
$ws.Range("I126").Value = 1311.7273
$ws.Range("J126").Value = 1413
$ws.Range("K126").Value = 3935.1819
$ws.Range("L126").Value = 4239
$ws.Range("M126").Value = -1465.1819
$ws.Range("N126").Value = -9179

$ws.Range("H132").Value = 61896.09
$ws.Range("I132").Value = 46227.816
$ws.Range("J132").Value = 93232.63
$ws.Range("K132").Value = 138683.448
$ws.Range("L132").Value = 279697.89
$ws.Range("M132").Value = -136153.448
$ws.Range("N132").Value = -284757.89

$ws.Range("H136").Value = 40353.49
$ws.Range("I136").Value = 25435
$ws.Range("J136").Value = 101519.3
$ws.Range("K136").Value = 76305
$ws.Range("L136").Value = 304557.9
$ws.Range("M136").Value = -73755
$ws.Range("N136").Value = -309657.9
